# Appends 5 new daily rows (2020-12-18 .. 2020-12-22) to the Indiana
# hospital ventilator dataset on the "Report" sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newRows = @(
    @{ Date = "2020-12-18"; Values = @(2150, 673, 1004, 473, 2806, 358, 446, 2002, 31.3, 46.7, 22, 12.76, 15.89, 71.34999999999999) },
    @{ Date = "2020-12-19"; Values = @(2154, 660, 986, 508, 2804, 343, 468, 1993, 30.64, 45.78, 23.58, 12.23, 16.69, 71.08) },
    @{ Date = "2020-12-20"; Values = @(2159, 665, 954, 540, 2806, 345, 450, 2011, 30.8, 44.19, 25.01, 12.3, 16.04, 71.67) },
    @{ Date = "2020-12-21"; Values = @(2145, 649, 1001, 495, 2809, 356, 487, 1966, 30.26, 46.67, 23.08, 12.67, 17.34, 69.98999999999999) },
    @{ Date = "2020-12-22"; Values = @(2125, 664, 1015, 446, 2807, 364, 466, 1977, 31.25, 47.76, 20.99, 12.97, 16.6, 70.43000000000001) }
)

$startRow = 297

for ($i = 0; $i -lt $newRows.Count; $i++) {
    $r = $startRow + $i
    $row = $newRows[$i]

    # Column A holds a plain text date label (matches the existing column,
    # which is stored as text, not a real date). Force text formatting
    # first so Excel does not auto-convert the string into a date serial,
    # then reset the style back to Normal so no stray number format is
    # left behind on the cell.
    $cellA = $ws.Cells.Item($r, 1)
    $cellA.NumberFormat = "@"
    $cellA.Value = $row.Date
    $cellA.Style = "Normal"

    for ($c = 0; $c -lt $row.Values.Count; $c++) {
        $ws.Cells.Item($r, $c + 2).Value = $row.Values[$c]
    }
}
